$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting the existing rows 7-16 down to 8-17.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44484
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 300000000
$ws.Range("G7").Value = "Espárragos"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 550
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 1700
$ws.Range("M7").Value = 1700
$ws.Range("N7").Value = "$/kilo"
$ws.Range("O7").Value = "Provincia de Linares"
$ws.Range("P7").Value = 1700
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
